$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.191.25'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').Value = '4.038.79'
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = "'539.32"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('D6').Value = "'149.25"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.26%  '
$ws.Range('D7').Value = '4.033.37'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +0.39%  '
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('D10').Value = "'0.751"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.42%  '
$ws.Range('E11').Value = '  -0.78%  '
$ws.Range('D12').Value = "'53.45"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +11.41%  '
$ws.Range('D13').Value = "'0.0000333"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.61%  '
$ws.Range('D14').Value = "'10.85"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('D15').Value = '4.676.89'
$ws.Range('E15').Value = '  -0.35%  '
$ws.Range('D16').Value = '4.043.45'
$ws.Range('E16').Value = '  +0.69%  '
$ws.Range('D17').Value = "'14.31"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.45%  '
$ws.Range('D18').Value = "'20.61"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.55%  '
$ws.Range('E19').Value = '  -1.04%  '
$ws.Range('E20').Value = '  -0.95%  '
$ws.Range('D21').Value = '72.098.08'
$ws.Range('E21').Value = '  +0.28%  '
$ws.Range('D22').Value = "'440.24"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.47%  '
$ws.Range('D23').Value = "'97.79"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.36%  '
$ws.Range('D24').Value = "'3.51"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.90%  '
$ws.Range('E25').Value = '  -0.31%  '
$ws.Range('D26').Value = "'14.58"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.19%  '
$ws.Range('E27').Value = '  +27.32%  '
$ws.Range('E28').Value = '  -0.76%  '
$ws.Range('D29').Value = "'10.70"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.61%  '
$ws.Range('E30').Value = '  +2.01%  '
$ws.Range('D31').Value = "'37.15"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.11%  '
$ws.Range('D32').Value = "'8.27"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +19.94%  '
$ws.Range('E33').Value = '  +1.67%  '
$ws.Range('D34').Value = "'13.52"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.93%  '
$ws.Range('D35').Value = "'49.41"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +14.94%  '
$ws.Range('D36').Value = "'683.43"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.78%  '
$ws.Range('D37').Value = "'66.74"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.68%  '
$ws.Range('E38').Value = '  +4.47%  '
$ws.Range('D39').Value = '0.0₃0914'
$ws.Range('E39').Value = '  +8.51%  '
$ws.Range('D40').Value = "'11.33"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +18.21%  '
$ws.Range('E41').Value = '  -6.58%  '
$ws.Range('B42').Value = 'ThetaToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('E42').Value = '  -1.31%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('E43').Value = '  +2.28%  '
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').Value = "'0.999"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('D46').Value = "'0.0492"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.91%  '
$ws.Range('E47').Value = '  -1.23%  '
$ws.Range('D48').Value = "'2.64"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.95%  '
$ws.Range('E49').Value = '  +2.32%  '
$ws.Range('B50').Value = 'FLOKI'
$ws.Range('C50').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D50').Value = "'0.000288"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.87%  '
$ws.Range('B51').Value = 'ApeXProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D51').Value = "'3.33"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.58%  '
